# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Leading apostrophe forces Excel to store the value as text,
    # preventing numeric-looking strings (e.g. "217.17") from being
    # auto-converted into numbers while keeping General number format.
    $ws.Range($cellRef).Value = "'" + $text
}

$ws.Range('D2').Value = '26.261.42'
$ws.Range('D3').Value = '1.675.91'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.25%  '
Set-TextCell 'D5' "217.17"
$ws.Range('E5').Value = '  +0.23%  '
Set-TextCell 'D6' "0.5334"
$ws.Range('E6').Value = '  +4.71%  '
$ws.Range('E7').Value = '  +0.24%  '
Set-TextCell 'D8' "0.2679"
$ws.Range('E8').Value = '  +1.45%  '
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('E10').Value = '  -0.17%  '
Set-TextCell 'D11' "0.07509"
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').Value = '1.682.05'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('E13').Value = '  +0.27%  '
Set-TextCell 'D14' "0.5761"
$ws.Range('E14').Value = '  -1.08%  '
Set-TextCell 'D15' "0.000008479"
$ws.Range('E15').Value = '  -0.18%  '
Set-TextCell 'D16' "64.57"
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').Value = '26.292.34'
$ws.Range('E17').Value = '  +0.87%  '
Set-TextCell 'D18' "4.908"
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  +1.13%  '
Set-TextCell 'D21' "190.11"
$ws.Range('E21').Value = '  +0.14%  '
Set-TextCell 'D22' "6.179"
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('E23').Value = '  +0.19%  '
Set-TextCell 'D24' "144.89"
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D25' "0.1277"
$ws.Range('E25').Value = '  +6.95%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D26' "7.799"
$ws.Range('E26').Value = '  +2.58%  '
Set-TextCell 'D27' "15.75"
$ws.Range('E27').Value = '  +0.80%  '
Set-TextCell 'D28' "0.06480"
$ws.Range('E28').Value = '  -3.53%  '
Set-TextCell 'D29' "1.364"
$ws.Range('E30').Value = '  +0.36%  '
Set-TextCell 'D31' "3.580"
$ws.Range('E31').Value = '  +1.76%  '
Set-TextCell 'D32' "3.583"
$ws.Range('E32').Value = '  +2.13%  '
Set-TextCell 'D33' "1.653"
Set-TextCell 'D34' "1.029"
$ws.Range('E34').Value = '  +1.16%  '
Set-TextCell 'D35' "0.6188"
$ws.Range('E35').Value = '  +1.79%  '
$ws.Range('E36').Value = '  +1.47%  '
Set-TextCell 'D37' "2.726"
$ws.Range('E37').Value = '  +0.53%  '
Set-TextCell 'D38' "6.292"
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').Value = '1.113.63'
$ws.Range('E39').Value = '  +3.53%  '
Set-TextCell 'D40' "0.01620"
$ws.Range('E40').Value = '  +1.19%  '
Set-TextCell 'D41' "0.8729"
$ws.Range('E41').Value = '  +1.60%  '
Set-TextCell 'D42' "1.014"
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').Value = '1.828.01'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 'D45' "0.00000000108"
$ws.Range('E45').Value = '  -5.69%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D46' "56.89"
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D47' "8.174"
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell 'D48' "1.001"
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D49' "0.05258"
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D50' "0.4289"
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D51' "6.074"
$ws.Range('E51').Value = '  +2.04%  '
